$d = $word.ActiveDocument

# Locate the paragraph that ends the "Osman" section (the one containing
# "...verändert"), right before the "Team" heading. The new bullet item
# needs to be inserted immediately after it, carrying the same list
# formatting (Listenabsatz / numId 4) and the trailing "_GoBack" bookmark.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*GUI Layout individuell angepasst und ver*ndert*") {
        $target = $p
    }
}

if ($target -eq $null) {
    throw "Could not find the 'GUI Layout ... verändert' paragraph"
}

# Insert a brand-new paragraph right after it; Word clones the
# paragraph/run formatting (pStyle Listenabsatz, numPr ilvl0/numId4,
# bold+underline rPr) from the split point, matching the target markup.
$target.Range.InsertParagraphAfter() | Out-Null

# Re-resolve the freshly inserted (still empty) paragraph as the one
# immediately following $target.
$newPara = $target.Next()

# Seed the paragraph with its final text plus one throw-away trailing
# character. We need a non-degenerate (Start != End) Range to reliably
# drive Bookmarks.Add, so we bookmark that placeholder character first
# and then shrink it away - this leaves the bookmark collapsed exactly
# where we want it (right after the real text) without ever handing a
# zero-length Range straight to Bookmarks.Add.
$newPara.Range.Text = "Logout Button für Kunde und Admin eingefügt#"

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$placeholderStart = $newPara.Range.End - 2
$placeholderEnd = $newPara.Range.End - 1
$bmRange = $d.Range($placeholderStart, $placeholderEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)

$bm = $d.Bookmarks("_GoBack")
$bm.Range.Text = ""
